$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-02 Sunday" "2025-11-03 Monday"

Replace-Text "688÷8=86, 0" "529÷8=66, 1"
Replace-Text "268÷8=33, 4" "248÷4=62, 0"
Replace-Text "307÷4=76, 3" "476÷6=79, 2"
Replace-Text "628÷4=157, 0" "787÷5=157, 2"
Replace-Text "646÷8=80, 6" "987÷5=197, 2"

Replace-Text "786÷2=393, 0" "545÷6=90, 5"
Replace-Text "960÷7=137, 1" "119÷3=39, 2"
Replace-Text "982÷9=109, 1" "550÷2=275, 0"
Replace-Text "918÷4=229, 2" "633÷5=126, 3"
Replace-Text "219÷9=24, 3" "995÷5=199, 0"

Replace-Text "544÷8=68, 0" "790÷5=158, 0"
Replace-Text "272÷9=30, 2" "513÷9=57, 0"
Replace-Text "222÷5=44, 2" "936÷5=187, 1"
Replace-Text "674÷4=168, 2" "632÷5=126, 2"
Replace-Text "425÷5=85, 0" "923÷5=184, 3"

Replace-Text "568÷8=71, 0" "441÷5=88, 1"
Replace-Text "617÷8=77, 1" "542÷9=60, 2"
Replace-Text "854÷9=94, 8" "581÷7=83, 0"
Replace-Text "249÷6=41, 3" "868÷4=217, 0"
Replace-Text "994÷8=124, 2" "946÷5=189, 1"

Replace-Text "560÷7=80, 0" "596÷3=198, 2"
Replace-Text "635÷8=79, 3" "613÷4=153, 1"
Replace-Text "209÷3=69, 2" "283÷3=94, 1"
Replace-Text "520÷6=86, 4" "663÷5=132, 3"
Replace-Text "101÷3=33, 2" "833÷8=104, 1"
